$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the month label (MAR-2015 -> FEB-2015) across the whole column
$ws.Range("A2:A280").Value = "'FEB-2015"

# Update selection / scroll position
$ws.Range("C5").Select()
$excel.ActiveWindow.ScrollRow = 1

# Autofit column C so its width matches its (now wider) content
$ws.Columns("C").AutoFit() | Out-Null
